$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = 32.50156583216571; E = 32.94283676147461; F = 34.21048527613666; G = 32.26889495557401 }
    3 = @{ D = 31.12023594651493; E = 32.72195434570312; F = 33.15878447514397; G = 30.08478366902539 }
    4 = @{ D = 33.00644847681398; E = 33.45501708984375; F = 34.60498068382673; G = 32.5334139556897 }
    5 = @{ D = 30.68970386348994; E = 30.76373291015625; F = 30.92824425297242; G = 27.93412753988001 }
    6 = @{ D = 33.22819603223149; E = 34.14981842041016; F = 34.87217077329465; G = 32.80474873228826 }
    7 = @{ D = 33.837934723348; E = 35.25273513793945; F = 36.56707016440913; G = 33.38587123025538 }
    8 = @{ D = 34.77268070374505; E = 33.26192474365234; F = 35.23687973956959; G = 33.11000489008421 }
    9 = @{ D = 37.33156166601932; E = 38.15699005126953; F = 38.71862421757658; G = 37.16136741038282 }
    10 = @{ D = 36.96005430659647; E = 38.95812606811523; F = 39.44692772737062; G = 36.41980315713099 }
    11 = @{ D = 40.00667988757762; E = 39.66097640991211; F = 40.23138746555141; G = 39.23748758607901 }
    12 = @{ D = 40.00943329659836; E = 40.78437042236328; F = 42.32553801833105; G = 39.84399724504603 }
    13 = @{ D = 38.34850301876148; E = 38.00625991821289; F = 42.20090341401771; G = 37.64647027034376 }
    14 = @{ D = 34.17838512874127; E = 35.60359191894531; F = 35.81604550902629; G = 34.03675156995863 }
    15 = @{ D = 35.83696122960632; E = 36.98031234741211; F = 38.03434100491589; G = 35.72083973424731 }
    16 = @{ D = 40.93739728803269; E = 41.62519454956055; F = 42.47377622574068; G = 39.65112525219609 }
    18 = @{ D = 53.34844018674714; E = 55.60538864135742; F = 55.65034507609375; G = 52.36833045205034 }
    19 = @{ D = 58.52706907906894; E = 60.73324584960938; F = 64.13293287788794; G = 58.201566281504 }
    20 = @{ D = 60.95915700445646; E = 61.14096450805664; F = 61.79546350658028; G = 59.16837780517545 }
    21 = @{ D = 65.42440096725305; E = 73.40030670166016; F = 74.27738501436652; G = 65.37871588952484 }
    22 = @{ D = 60.8683261867532; E = 64.73342132568359; F = 66.97351921304117; G = 58.47215134430768 }
    23 = @{ D = 63.05701663092425; E = 64.47879791259766; F = 65.54975656568148; G = 61.57983766155381 }
    24 = @{ D = 74.38973021233218; E = 72.23699951171875; F = 77.71161420514916; G = 71.54107460838905 }
    25 = @{ D = 76.43759578898127; E = 76.02738952636719; F = 77.62158964312036; G = 72.75508907059996 }
    26 = @{ D = 82.23480186773816; E = 86.83358001708984; F = 87.28315635977556; G = 82.11304218171169 }
    27 = @{ D = 93.5946593959798; E = 99.69947052001952; F = 99.9816696196051; G = 92.51290654748394 }
    28 = @{ D = 104.7134850100623; E = 113.5867080688477; F = 113.983168792034; G = 102.193114180154 }
    29 = @{ D = 112.8881983413681; E = 101.6836700439453; F = 112.9734418216001; G = 96.45552387300648 }
    30 = @{ D = 103.3973992194267; E = 95.10124969482422; F = 104.5390737374864; G = 94.84437467575744 }
    31 = @{ D = 98.9728622916322; E = 107.3169784545898; F = 108.0720564916435; G = 96.83187459488266 }
    32 = @{ D = 99.14982521372704; E = 101.9706192016602; F = 103.2370926682909; G = 93.6329728413386 }
    33 = @{ D = 109.5204116287298; E = 114.6067657470703; F = 115.4737568003429; G = 109.5204116287298 }
    34 = @{ D = 125.6979779583206; E = 127.197135925293; F = 129.1799106580013; G = 124.7404372279016 }
    35 = @{ D = 130.1361101824991; E = 131.1349182128906; F = 131.193098906433; G = 125.2390329208732 }
    36 = @{ D = 118.359173386968; E = 116.6850128173828; F = 124.9195400464748; G = 114.3295112839243 }
    37 = @{ D = 141.419740521661; E = 145.190673828125; F = 148.9420793559621; G = 139.504974747138 }
    38 = @{ D = 152.5038588308085; E = 148.1124420166016; F = 152.7587127008474; G = 146.5440713901251 }
    39 = @{ D = 165.4261760174149; E = 159.2308044433594; F = 167.4617893171745; G = 157.4213670646024 }
    40 = @{ D = 168.7096486664738; E = 161.4482269287109; F = 168.9365638617247; G = 161.3495602936897 }
    41 = @{ D = 158.7058635782533; E = 166.5560150146484; F = 168.8328610343947; G = 150.0736576094974 }
    42 = @{ D = 171.6419996777651; E = 184.744369506836; F = 185.2013199976124; G = 158.7085137299506 }
    43 = @{ D = 199.7503963180078; E = 200.2187805175781; F = 217.6685325614052; G = 196.3321891938854 }
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item([int]$r, 4).Value = $vals.D
    $ws.Cells.Item([int]$r, 5).Value = $vals.E
    $ws.Cells.Item([int]$r, 6).Value = $vals.F
    $ws.Cells.Item([int]$r, 7).Value = $vals.G
    $ws.Cells.Item([int]$r, 8).Value = 450687724
    $ws.Cells.Item([int]$r, 9).Value = "TRI"
}
